# Update the "Progress" value for the OPERA-2 trial (row 10, column B)
# in the clinical trials setup sheet from 12 to 25.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("B10").Value = 25
